$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G3").Value = "2016-08-18 00:43:27"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H3").Value = "2016-08-18 00:43:22"
$zhcn.Range("K3").Value = "2016-08-18 00:43:40"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H3").Value = "2016-08-18 00:43:27"
$dede.Range("K3").Value = "2016-08-18 00:43:47"
